# "Add_Products" sheet gains a new trailing identifier column (AE), mirroring
# the existing SKU identifier columns (AA/AC/AD), and the product identifier
# value previously stored in AB2 is replaced by a freshly generated one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add_Products")

# New column AE: header matches the other "SKU " identifier headers (D1, AA1,
# AC1, AD1) and the data row gets a newly generated identifier value.
$ws.Range("AE1").Value = "SKU "
$ws.Range("AE2").Value = "EcomvMFj"

# Existing identifier value in AB2 is regenerated.
$ws.Range("AB2").Value = "ProdId0Svp"
